$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "603.42")
# are stored as text, matching the inline-string cells in the workbook.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '68.001.88'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '3.531.58'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '603.42'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').Value = '182.10'
$ws.Range('E6').Value = '  +4.81%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.531.30'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  +5.74%  '
$ws.Range('D11').Value = '7.18'
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').Value = '0.441'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '4.145.22'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '32.43'
$ws.Range('E14').Value = '  +10.48%  '
$ws.Range('D15').Value = '0.136'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '67.992.45'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').Value = '0.0000181'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').Value = '3.538.49'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').Value = '6.40'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').Value = '14.57'
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('D21').Value = '401.48'
$ws.Range('E21').Value = '  +1.35%  '
$ws.Range('D22').Value = '8.06'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = '74.06'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = '0.548'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = '5.72'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').Value = '0.0000124'
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('D28').Value = '10.60'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('D29').Value = '0.178'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '6.34'
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('D32').Value = '1.46'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').Value = '2.10'
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('D34').Value = '24.06'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '7.55'
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = '1.65'
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('D38').Value = '163.16'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '0.882'
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').Value = '2.82'
$ws.Range('E41').Value = '  +7.84%  '
$ws.Range('D42').Value = '7.03'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '4.76'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').Value = '2.906.33'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').Value = '26.62'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '0.0741'
$ws.Range('E46').Value = '  -1.58%  '
$ws.Range('D47').Value = '27.02'
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('D48').Value = '42.61'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('D49').Value = '352.06'
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('D50').Value = '0.0307'
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').Value = '1.08'
$ws.Range('E51').Value = '  -0.93%  '

# Restore the default "Normal" style on column D now that the values are
# already stored as text, so no stray number-format/style is left behind.
$dRange.Style = "Normal"
